$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VOC")

# --- Insert a new column before D (a new, most-recent fiscal-year column). ---
# Excel shifts existing D:K data right to E:L automatically and preserves
# row/formula structure; we still need to carry the formatting across and
# populate the new column with the new period's figures.
$ws.Columns("D").Insert()

# Pull the number/style formatting from column E (which now holds what used
# to be column D) into the freshly inserted, blank column D.
$ws.Columns("E").Copy()
$ws.Columns("D").PasteSpecial(-4122)   # xlPasteFormats

# Seed column D with column E's values (this correctly carries forward every
# row whose figure is unchanged/constant across the table, e.g. the "NA"
# placeholders, the 0 placeholders, and the blank spacer rows).
$ws.Range("D7:D102").Value2 = $ws.Range("E7:E102").Value2

$excel.CutCopyMode = 0

# --- Now overwrite column D with the genuinely new figures for this period ---

# "Period Ending" header rows (new date 2018-12-31)
$ws.Range("D7").Value2 = 43465
$ws.Range("D38").Value2 = 43465
$ws.Range("D80").Value2 = 43465

# Income statement block
$ws.Range("D8").Value2 = 12900    # Total Revenue
$ws.Range("D17").Value2 = 900     # Research Development
$ws.Range("D18").Value2 = 12100   # Selling General and Administrative
$ws.Range("D21").Value2 = 12100   # Total Operating Expenses
$ws.Range("D23").Value2 = 12100   # Operating Income or Loss
$ws.Range("D26").Value2 = 12100   # Earnings Before Interest And Taxes
$ws.Range("D27").Value2 = 12100   # Income Before Tax
$ws.Range("D33").Value2 = 12100   # Net Income
$ws.Range("D35").Value2 = 12100   # Net Income Applicable To Common Shares

# Balance sheet block
$ws.Range("D41").Value2 = 400     # Cash And Cash Equivalents (first numeric balance-sheet row)
$ws.Range("D49").Value2 = 69700   # Total Current Assets
$ws.Range("D54").Value2 = 70000   # Total Assets
$ws.Range("D76").Value2 = 70000   # Total liabilities and stockholders' equity

# Cash flow block
$ws.Range("D81").Value2 = 12100   # Net Income (restated at top of cash-flow statement)
